# Update cryptocurrency Price (column D) and Volume(1h) (column E) values
# for the latest scrape, per GitHub Actions automation run.
# Values are stored as plain text (inlineStr) in the workbook, so we
# prefix assignments with an apostrophe to force text entry (preventing
# Excel from auto-converting "41.03" to a number or "0.69%" to a percentage),
# then reset the cell style to Normal so no extra number formatting sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.69%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.46%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.112"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.84%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07637"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'1.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.29%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'2.484"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.04%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9051"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.31%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1112"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'7.58%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1784"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.37%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09075"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.47%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04246"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.46%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1050"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.47%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.59%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005749"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.60%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.340"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.43%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.253"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.60%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.3318"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.24%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.643"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-4.90%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1360"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.94%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2709"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.04%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04038"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.46%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001256"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'4.64%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004111"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.60%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001300"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.04%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D38").Value = "'0.02400"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'0.83%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05231"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.24%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007792"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.32%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1300"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.37%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007043"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'20.60%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001949"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.01%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008427"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'6.05%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3329"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.01%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006911"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'7.77%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05577"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1,270.25%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.05%"
$ws.Range("E51").Style = "Normal"
